$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1508.2727
$ws.Range("I42").Value = 188.875
$ws.Range("J42").Value = 5026.6665
$ws.Range("K42").Value = 566.625
$ws.Range("L42").Value = 15079.9995
$ws.Range("M42").Value = -336.625
$ws.Range("N42").Value = -15539.9995

$ws.Range("H62").Value = 4512.846
$ws.Range("I62").Value = 4664
$ws.Range("K62").Value = 4664
$ws.Range("M62").Value = -4040

$ws.Range("H65").Value = 4512.846
$ws.Range("I65").Value = 4664
$ws.Range("K65").Value = 23320
$ws.Range("M65").Value = -20200

$ws.Range("H127").Value = 594.4
$ws.Range("I127").Value = 478.5
$ws.Range("K127").Value = 1435.5
$ws.Range("M127").Value = 3524.5

$ws.Range("H135").Value = 662.9
$ws.Range("I135").Value = 561.7895
$ws.Range("K135").Value = 5056.1055
$ws.Range("M135").Value = -2521.1055

$ws.Range("H138").Value = 3952.532
$ws.Range("J138").Value = 4547.4614
$ws.Range("L138").Value = 13642.3842
$ws.Range("N138").Value = -23922.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4101.375
$ws.Range("I32").Value = 3437.3333
$ws.Range("K32").Value = 3437.3333
$ws.Range("M32").Value = -3150.3333

$ws.Range("H61").Value = 4961.6
$ws.Range("I61").Value = 3969.6667
$ws.Range("K61").Value = 3969.6667
$ws.Range("M61").Value = -3757.6667

$ws.Range("H102").Value = 3023.1538
$ws.Range("I102").Value = 3317.6365
$ws.Range("K102").Value = 3317.6365
$ws.Range("M102").Value = -1695.6365

$ws.Range("H107").Value = 69942.5
$ws.Range("J107").Value = 69890
$ws.Range("L107").Value = 69890
$ws.Range("N107").Value = -77570

$ws.Range("H109").Value = 69890
$ws.Range("J109").Value = 69890
$ws.Range("L109").Value = 69890
$ws.Range("N109").Value = -72664

$ws.Range("H117").Value = 30160.25
$ws.Range("J117").Value = 30160.25
$ws.Range("L117").Value = 30160.25
$ws.Range("N117").Value = -39338.25

$ws.Range("H122").Value = 5939.4136
$ws.Range("I122").Value = 6015.8213
$ws.Range("K122").Value = 18047.4639
$ws.Range("M122").Value = -15597.4639

$ws.Range("H132").Value = 8526.833000000001
$ws.Range("I132").Value = 6975.273
$ws.Range("J132").Value = 11940.267
$ws.Range("K132").Value = 20925.819
$ws.Range("L132").Value = 35820.801
$ws.Range("M132").Value = -18395.819
$ws.Range("N132").Value = -40880.801

$ws.Range("H136").Value = 4961.6
$ws.Range("I136").Value = 3969.6667
$ws.Range("K136").Value = 11909.0001
$ws.Range("M136").Value = -9359.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 9247.333000000001
$ws.Range("J106").Value = 9247.333000000001
$ws.Range("L106").Value = 9247.333000000001
$ws.Range("N106").Value = -11771.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 100000
$ws.Range("J75").Value = 100000
$ws.Range("L75").Value = 100000
$ws.Range("N75").Value = -101996

$ws.Range("H78").Value = 100000
$ws.Range("J78").Value = 100000
$ws.Range("L78").Value = 300000
$ws.Range("N78").Value = -309984

$ws.Range("H93").Value = 15849
$ws.Range("I93").Value = 14465.333
$ws.Range("K93").Value = 14465.333
$ws.Range("M93").Value = -12593.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7505.7334
$ws.Range("J39").Value = 8553.5
$ws.Range("L39").Value = 25660.5
$ws.Range("N39").Value = -26248.5

$ws.Range("H63").Value = 12116.5
$ws.Range("I63").Value = 12159.667
$ws.Range("K63").Value = 36479.001
$ws.Range("M63").Value = -35730.001

$ws.Range("H66").Value = 12116.5
$ws.Range("I66").Value = 12159.667
$ws.Range("K66").Value = 109437.003
$ws.Range("M66").Value = -105693.003

$ws.Range("H86").Value = 633.7
$ws.Range("I86").Value = 479.4
$ws.Range("K86").Value = 1438.2
$ws.Range("M86").Value = -252.1999999999998

$ws.Range("H89").Value = 633.7
$ws.Range("I89").Value = 479.4
$ws.Range("K89").Value = 4314.599999999999
$ws.Range("M89").Value = 1613.400000000001

$ws.Range("H137").Value = 2968.2144
$ws.Range("J137").Value = 3296.875
$ws.Range("L137").Value = 9890.625
$ws.Range("N137").Value = -20090.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 386.44446
$ws.Range("I97").Value = 417.2
$ws.Range("J97").Value = 298.57144
$ws.Range("K97").Value = 417.2
$ws.Range("L97").Value = 298.57144
$ws.Range("M97").Value = 78.80000000000001
$ws.Range("N97").Value = -1290.57144

$ws.Range("H132").Value = 6441.727
$ws.Range("I132").Value = 3837
$ws.Range("K132").Value = 11511
$ws.Range("M132").Value = -8981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 74983.664
$ws.Range("I5").Value = 74983.664
$ws.Range("K5").Value = 74983.664
$ws.Range("M5").Value = -74870.664

$ws.Range("H40").Value = 5293.8623
$ws.Range("I40").Value = 5397
$ws.Range("K40").Value = 5397
$ws.Range("M40").Value = -5261

$ws.Range("H64").Value = 77163.57000000001
$ws.Range("J64").Value = 77163.57000000001
$ws.Range("L64").Value = 77163.57000000001
$ws.Range("N64").Value = -77613.57000000001

$ws.Range("H67").Value = 77163.57000000001
$ws.Range("J67").Value = 77163.57000000001
$ws.Range("L67").Value = 77163.57000000001
$ws.Range("N67").Value = -78723.57000000001

$ws.Range("H74").Value = 39992.25
$ws.Range("I74").Value = 39984
$ws.Range("K74").Value = 39984
$ws.Range("M74").Value = -38986

$ws.Range("H77").Value = 39992.25
$ws.Range("I77").Value = 39984
$ws.Range("K77").Value = 119952
$ws.Range("M77").Value = -114960

$ws.Range("H104").Value = 99000
$ws.Range("J104").Value = 99000
$ws.Range("L104").Value = 99000
$ws.Range("N104").Value = -105988

$ws.Range("H106").Value = 39798
$ws.Range("J106").Value = 39798
$ws.Range("L106").Value = 39798
$ws.Range("N106").Value = -42322

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -1280

$ws.Range("H13").Value = 17520.8
$ws.Range("I13").Value = 24534.666
$ws.Range("J13").Value = 7000
$ws.Range("K13").Value = 24534.666
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = -24394.666
$ws.Range("N13").Value = -7280
